# Generate Report for Archive
# Update the localization status for file cad4860f-87f8-48a8-8617-6517a6f51269.md
# from "Ready for handoff" to "In Translation" across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E7").Value = "In Translation"
$wsOverview.Range("F7").Value = "In Translation"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C7").Value = "In Translation"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C7").Value = "In Translation"
